$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Update header labels: append unit/precision hints to the
#    column headers (e.g. "电压（峰峰值）/V" -> "电压（峰峰值）/V 2dp")
# -----------------------------------------------------------------
$ws.Range("B8").Value  = "电压（峰峰值）/V 2dp"
$ws.Range("C8").Value  = "周期/kHz 3sd"
$ws.Range("D8").Value  = "频率/ms 3dp"

$ws.Range("B13").Value = "电压（峰峰值）/V 2dp"
$ws.Range("C13").Value = "周期/Hz 3sd"
$ws.Range("D13").Value = "频率/ms 3dp"

$ws.Range("B18").Value = "函数发生器频率f信/Hz 2dp"
$ws.Range("C18").Value = "算出的市电频率fx/Hz 2dp"
$ws.Range("D18").Value = "平均市电频率/Hz 2dp"

# -----------------------------------------------------------------
# 2. Apply number formats to the data cells so the displayed
#    precision matches the new header hints.
#    - voltage / averaged-frequency columns -> 2 decimal places
#    - raw period (ms) columns              -> 3 decimal places
#    - derived kHz/Hz columns                -> scientific (3 sig figs)
# -----------------------------------------------------------------
$cells2dp = @("B9","B10","B11","B14","B15","B16","B19","C19","D19","B20","C20","B21","C21","B22","C22")
foreach ($cellAddr in $cells2dp) {
    $ws.Range($cellAddr).NumberFormat = "0.00_ "
}

$cells3dp = @("D9","D10","D11","D14","D15","D16")
foreach ($cellAddr in $cells3dp) {
    $ws.Range($cellAddr).NumberFormat = "0.000_ "
}

$cellsSci = @("C9","C10","C11","C14","C15","C16")
foreach ($cellAddr in $cellsSci) {
    $ws.Range($cellAddr).NumberFormat = "0.00E+00"
}

# -----------------------------------------------------------------
# 3. Widen columns B and C (no longer auto "best fit") to make
#    room for the longer header text.
# -----------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 26.21875
$ws.Columns("C").ColumnWidth = 24.88671875

# -----------------------------------------------------------------
# 4. Update the window view: scroll so row 3 is at the top and
#    select cell E21.
# -----------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("E21").Select() | Out-Null
